$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the "last changed" date for every data row
# (rows 2-236). This whole column bumps from 2023-10-04 (45203) to
# 2023-10-05 (45204) because the source feed was refreshed a day later.
for ($r = 2; $r -le 236; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cell.Value = 45204
}

# Row 236 picks up an explicit (custom) row height of 15 in the refreshed
# export - same visual height as before, just now stored explicitly.
$ws.Rows.Item(236).RowHeight = 15

# A new case was appended as row 237: "A 47362-2023".
$ws.Cells.Item(237, 1).Value = "A 47362-2023"

$ws.Cells.Item(237, 2).Value = 45202
$ws.Cells.Item(237, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(237, 3).Value = 45204
$ws.Cells.Item(237, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item(237, 4).Value = "HALLANDS LÄN"
$ws.Cells.Item(237, 5).Value = "KUNGSBACKA"

$ws.Cells.Item(237, 7).Value = 10.8
$ws.Cells.Item(237, 8).Value = 0
$ws.Cells.Item(237, 9).Value = 0
$ws.Cells.Item(237, 10).Value = 0
$ws.Cells.Item(237, 11).Value = 0
$ws.Cells.Item(237, 12).Value = 0
$ws.Cells.Item(237, 13).Value = 0
$ws.Cells.Item(237, 14).Value = 0
$ws.Cells.Item(237, 15).Value = 0
$ws.Cells.Item(237, 16).Value = 0
$ws.Cells.Item(237, 17).Value = 0

# R237 mirrors the empty, wrap-text-styled "Artnamn" cell seen on every
# other data row (e.g. R236) that has no species listed.
$ws.Cells.Item(237, 18).Value = ""
$ws.Cells.Item(237, 18).WrapText = $true
